$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Slide 1: remove the "Cloud 50" shape and its "Elbow Connector 51"
#    connector (the small unused "Web" cloud callout wired to shape id 2).
# ---------------------------------------------------------------------------
$s = $p.Slides.Item(1)

$connector = $s.Shapes.Item("Elbow Connector 51")
$connector.Delete()

$cloud = $s.Shapes.Item("Cloud 50")
$cloud.Delete()

# ---------------------------------------------------------------------------
# 2. Refresh the "datetimeFigureOut" date placeholders (slide master, every
#    slide layout, and the notes master) from 10/16/2016 to 3/27/17 - this
#    mirrors PowerPoint re-stamping the auto date field on every layout when
#    the deck was next saved.
# ---------------------------------------------------------------------------
$newDate = "3/27/17"

$master = $p.SlideMaster
$master.Shapes.Item("Date Placeholder 3").TextFrame.TextRange.Text = $newDate

$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $layout = $layouts.Item($i)
    for ($j = 1; $j -le $layout.Shapes.Count; $j++) {
        $shape = $layout.Shapes.Item($j)
        if ($shape.HasTextFrame -and $shape.TextFrame.TextRange.Text -match "10/16/2016") {
            $shape.TextFrame.TextRange.Text = $newDate
        }
    }
}

$notesMaster = $p.NotesMaster
$notesMaster.Shapes.Item("Date Placeholder 2").TextFrame.TextRange.Text = $newDate
